$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Create the new "2022-Q4" sheet by copying the existing "2022-Q3" sheet
#    (placed immediately before it). Copying - rather than adding a blank
#    sheet - carries over the correct sheetPr / pageMargins / column
#    formatting and (crucially) the existing text-typed cells, so the new
#    data can be dropped straight into already-text cells without Excel's
#    automatic "looks like a number" reinterpretation.
# ---------------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($q3)
$q4 = $wb.Worksheets.Item("2022-Q3 (2)")
$q4.Name = "2022-Q4"

# The copied sheet has 3 data rows (like 2022-Q3); 2022-Q4 only needs 2, so
# drop the extra row.
$q4.Rows("4:4").Delete()

# ---------------------------------------------------------------------------
# 2. Overwrite the copied sheet's data with the 2022-Q4 fund holdings.
#    Columns B-G are text in the source data (codes / formatted numbers
#    with fixed decimals, e.g. "90.00"), so force text formatting before
#    writing, then drop the formatting override back to Normal so no
#    residual number-format style is left applied to the cells.
# ---------------------------------------------------------------------------
$textCells = $q4.Range("B2:G3")
$textCells.NumberFormat = "@"

$q4.Range("A2").Value = 0
$q4.Range("B2").Value = "005638"
$q4.Range("C2").Value = "农银汇理量化智慧动力混合"
$q4.Range("D2").Value = "0.46"
$q4.Range("E2").Value = "90.00"
$q4.Range("F2").Value = "2.65"
$q4.Range("G2").Value = "0.0122"
$q4.Range("H2").Value = 2

$q4.Range("A3").Value = 1
$q4.Range("B3").Value = "562530"
$q4.Range("C3").Value = "华夏中证智选1000价值稳健策略ETF"
$q4.Range("D3").Value = "0.36"
$q4.Range("E3").Value = "96.22"
$q4.Range("F3").Value = "0.94"
$q4.Range("G3").Value = "0.0034"
$q4.Range("H3").Value = 6

$textCells.Style = "Normal"

# ---------------------------------------------------------------------------
# 3. Insert a new row into "总计" for the 2022-Q4 summary, pushing the
#    existing 2022-Q3 / 2022-Q2 rows down by one. Values are written
#    literally for every cell so the final content matches exactly. The
#    new A4 cell needs the same formatting as the existing A2/A3 label
#    cells, so its format is copied across explicitly.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

$total.Range("A3").Copy()
$total.Range("A4").PasteSpecial(-4122)

$total.Range("A4").Value = 2
$total.Range("B4").Value = "2022-Q2"
$total.Range("C4").Value = 2
$total.Range("D4").Value = 0.09

$total.Range("A3").Value = 1
$total.Range("B3").Value = "2022-Q3"
$total.Range("C3").Value = 3
$total.Range("D3").Value = 0.14

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 0.02

# ---------------------------------------------------------------------------
# 4. Restore the original active-sheet bookkeeping: the workbook's active
#    tab stays on "总计" and "2022-Q2" keeps being the sheet-view-selected
#    (tabSelected) sheet, exactly as in the source file - the new sheet
#    should not steal either of those just because it was created last.
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("2022-Q2").Activate()
$total.Activate()
